$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "type" column (D). This shifts registry (E->D) and tags (F->E) left.
$ws.Columns.Item(4).Delete()

# The registry column now holds namespaced identifiers instead of bare ids.
$ws.Cells.Item(4, 4).Value2 = "biotools:DAISY"
$ws.Cells.Item(16, 4).Value2 = "biotools:phyre"
$ws.Cells.Item(17, 4).Value2 = "fairsharing:nd9ce9"

# Refresh the autofilter so its range no longer includes the removed column.
$ws.AutoFilterMode = $false
$ws.Range("A1:E1").AutoFilter() | Out-Null

# Keep the _FilterDatabase defined name in sync with the new autofilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "tools!_FilterDatabase") {
        $n.RefersTo = "=tools!`$A`$1:`$E`$1"
    }
}

# Update the active selection to reflect where the author left the cursor.
$ws.Range("D17").Select() | Out-Null
